# Generate Report for Handback
# The handback transform failed for ab98151a-0aef-40f9-815a-efadcfc413ec:
# the handback file name (3p5xa3ss.axs) didn't match the handoff file name.
# Update the Status for that file on every sheet, and record the error
# detail on the per-locale (zh-cn / de-de) sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Status column for the ab98151a-0aef-40f9-815a-efadcfc413ec.md row (row 3)
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Error Detail column (K) for the same row, per locale
$wsZhCn.Range("K3").Value = "Handback file name: 3p5xa3ss.axs is different with handoff file name: ab98151a-0aef-40f9-815a-efadcfc413ec.bffb2ca696518a988bc792fdb6bb9f30da0a7c5d.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: 3p5xa3ss.axs is different with handoff file name: ab98151a-0aef-40f9-815a-efadcfc413ec.bffb2ca696518a988bc792fdb6bb9f30da0a7c5d.de-de."
